$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Collin Sexton"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Utah Jazz"

$ws.Range("A6").Value = "Anthony Edwards"
$ws.Range("B6").Value = "SG,SF"
$ws.Range("C6").Value = "Minnesota Timberwolves"

$ws.Range("A7").Value = "Jaden McDaniels"
$ws.Range("B7").Value = "SF,PF"
$ws.Range("C7").Value = "Minnesota Timberwolves"

$ws.Range("A10").Value = "Ivica Zubac"
$ws.Range("B10").Value = "C"
$ws.Range("C10").Value = "LA Clippers"

$ws.Range("A11").Value = "Wendell Carter Jr."
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Orlando Magic"

$ws.Range("A14").Value = "Jaren Jackson Jr."
$ws.Range("B14").Value = "PF,C"
$ws.Range("C14").Value = "Memphis Grizzlies"

$ws.Range("A15").Value = "Draymond Green"
$ws.Range("B15").Value = "PF,C"
$ws.Range("C15").Value = "Golden State Warriors"

$ws.Range("A16").Value = "Jayson Tatum"
$ws.Range("B16").Value = "SF,PF"
$ws.Range("C16").Value = "Boston Celtics"

$ws.Range("A17").Value = "Aaron Gordon"
$ws.Range("B17").Value = "PF,C"
$ws.Range("C17").Value = "Denver Nuggets"

$ws.Range("A18").Value = "Anfernee Simons"
$ws.Range("B18").Value = "PG,SG"
$ws.Range("C18").Value = "Portland Trail Blazers"

$ws.Range("A19").Value = "Paul George"
$ws.Range("B19").Value = "SG,SF,PF"
$ws.Range("C19").Value = "Philadelphia 76ers"
